$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.009582666666667
$ws.Range("H2").Value = 3.028748
$ws.Range("I2").Value = 0.2254436189979109
$ws.Range("J2").Value = 0.2254436189979109
$ws.Range("M2").Value = 0.6068319999999999
$ws.Range("N2").Value = 1.820496
$ws.Range("O2").Value = 0.03392274820144286
$ws.Range("P2").Value = 0.03392274820144286
$ws.Range("Q2").Value = 0.6126470687786665
$ws.Range("R2").Value = 5.513823619007999
$ws.Range("S2").Value = 0.007647667120888151
$ws.Range("T2").Value = 0.007647667120888152
$ws.Range("G3").Value = 1.009582666666667
$ws.Range("H3").Value = 3.028748
$ws.Range("I3").Value = 0.2254436189979109
$ws.Range("J3").Value = 0.2254436189979109
$ws.Range("O3").Value = 0.4504903529585388
$ws.Range("P3").Value = 0.4504903529585388
$ws.Range("Q3").Value = 8.135885471725333
$ws.Range("R3").Value = 73.22296924552801
$ws.Range("S3").Value = 0.1015601754946192
$ws.Range("T3").Value = 0.1015601754946192
$ws.Range("G4").Value = 1.009582666666667
$ws.Range("H4").Value = 3.028748
$ws.Range("I4").Value = 0.2254436189979109
$ws.Range("J4").Value = 0.2254436189979109
$ws.Range("M4").Value = 9.223151
$ws.Range("N4").Value = 27.669453
$ws.Range("O4").Value = 0.5155868988400183
$ws.Range("P4").Value = 0.5155868988400183
$ws.Range("Q4").Value = 9.311533381649332
$ws.Range("R4").Value = 83.80380043484399
$ws.Range("S4").Value = 0.1162357763824035
$ws.Range("T4").Value = 0.1162357763824035
$ws.Range("I5").Value = 0.4390905462561113
$ws.Range("J5").Value = 0.4390905462561113
$ws.Range("M5").Value = 0.6068319999999999
$ws.Range("N5").Value = 1.820496
$ws.Range("O5").Value = 0.03392274820144286
$ws.Range("P5").Value = 0.03392274820144286
$ws.Range("Q5").Value = 1.193236416661333
$ws.Range("R5").Value = 10.739127749952
$ws.Range("S5").Value = 0.01489515803828006
$ws.Range("T5").Value = 0.01489515803828006
$ws.Range("I6").Value = 0.4390905462561113
$ws.Range("J6").Value = 0.4390905462561113
$ws.Range("O6").Value = 0.4504903529585388
$ws.Range("P6").Value = 0.4504903529585388
$ws.Range("S6").Value = 0.1978060551636732
$ws.Range("T6").Value = 0.1978060551636732
$ws.Range("I7").Value = 0.4390905462561113
$ws.Range("J7").Value = 0.4390905462561113
$ws.Range("M7").Value = 9.223151
$ws.Range("N7").Value = 27.669453
$ws.Range("O7").Value = 0.5155868988400183
$ws.Range("P7").Value = 0.5155868988400183
$ws.Range("Q7").Value = 18.13582614227067
$ws.Range("R7").Value = 163.222435280436
$ws.Range("S7").Value = 0.226389333054158
$ws.Range("T7").Value = 0.226389333054158
$ws.Range("G8").Value = 1.502284666666667
$ws.Range("H8").Value = 4.506854000000001
$ws.Range("I8").Value = 0.3354658347459779
$ws.Range("J8").Value = 0.3354658347459779
$ws.Range("M8").Value = 0.6068319999999999
$ws.Range("N8").Value = 1.820496
$ws.Range("O8").Value = 0.03392274820144286
$ws.Range("P8").Value = 0.03392274820144286
$ws.Range("Q8").Value = 0.9116344088426667
$ws.Range("R8").Value = 8.204709679584001
$ws.Range("S8").Value = 0.01137992304227465
$ws.Range("T8").Value = 0.01137992304227465
$ws.Range("G9").Value = 1.502284666666667
$ws.Range("H9").Value = 4.506854000000001
$ws.Range("I9").Value = 0.3354658347459779
$ws.Range("J9").Value = 0.3354658347459779
$ws.Range("O9").Value = 0.4504903529585388
$ws.Range("P9").Value = 0.4504903529585388
$ws.Range("Q9").Value = 12.10640435644934
$ws.Range("R9").Value = 108.957639208044
$ws.Range("S9").Value = 0.1511241223002464
$ws.Range("T9").Value = 0.1511241223002464
$ws.Range("G10").Value = 1.502284666666667
$ws.Range("H10").Value = 4.506854000000001
$ws.Range("I10").Value = 0.3354658347459779
$ws.Range("J10").Value = 0.3354658347459779
$ws.Range("M10").Value = 9.223151
$ws.Range("N10").Value = 27.669453
$ws.Range("O10").Value = 0.5155868988400183
$ws.Range("P10").Value = 0.5155868988400183
$ws.Range("Q10").Value = 13.85579832565134
$ws.Range("R10").Value = 124.702184930862
$ws.Range("S10").Value = 0.1729617894034568
$ws.Range("T10").Value = 0.1729617894034568